# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gains a new (blank) column before the
# existing "Late" column, pushing "Late", "heading"(=Paid Date label reused
# for a blank-header column) and "Outstanding" one column to the right.
#
# Concretely: insert a new column N (shifting old N->O, O->P, P->Q) and give
# it the same width as column M, leaving the new column N's header/data cells
# blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the width Excel would assign to the freshly-inserted column (same as
# the neighbouring "Outstanding/header" column to its left).
$newColWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N, shifting the old N:P columns (and their
# widths/contents) one position to the right (N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $newColWidth

# Restore the selection left behind by the edit.
$ws.Range("R8").Select() | Out-Null
